# Updates the weekly Fruta/Granada price sheet (Macroferia Regional de Talca).
# The commit reshuffles the per-record values of Fecha/Calidad/Volumen/Precios/
# Unidad/Origen/Precio-Kg/Kg-unidad across the existing rows (row identity columns
# A, B, C, E-K are unchanged), so we only need to overwrite the affected cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("S2").Value = 944
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = (Get-Date -Year 2021 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 1111

# Row 4
$ws.Range("D4").Value = (Get-Date -Year 2021 -Month 5 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 1111

# Row 5
$ws.Range("D5").Value = (Get-Date -Year 2021 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M5").Value = 230
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 1111

# Row 6
$ws.Range("D6").Value = (Get-Date -Year 2021 -Month 6 -Day 7 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 18000
$ws.Range("P6").Value = 18000
$ws.Range("S6").Value = 1000

# Row 7
$ws.Range("D7").Value = (Get-Date -Year 2021 -Month 5 -Day 26 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L7").Value = 'Especial'
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 20000
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 20000
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1111
$ws.Range("T7").Value = 18

# Row 8
$ws.Range("D8").Value = (Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 300

# Row 9
$ws.Range("D9").Value = (Get-Date -Year 2021 -Month 4 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 800

# Row 10
$ws.Range("D10").Value = (Get-Date -Year 2021 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("R10").Value = 'Región Metropolitana'
$ws.Range("S10").Value = 800
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = (Get-Date -Year 2022 -Month 5 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 150
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 17000
$ws.Range("S11").Value = 944

# Row 12
$ws.Range("D12").Value = (Get-Date -Year 2022 -Month 6 -Day 2 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100

# Row 13
$ws.Range("D13").Value = (Get-Date -Year 2021 -Month 4 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("Q13").Value = '$/caja 15 kilos granel'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 1000
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = (Get-Date -Year 2021 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("D15").Value = (Get-Date -Year 2021 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("S15").Value = 944

# Row 16
$ws.Range("D16").Value = (Get-Date -Year 2021 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 20000
$ws.Range("P16").Value = 20000
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1111
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = (Get-Date -Year 2022 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("Q17").Value = '$/caja 15 kilos granel'
$ws.Range("T17").Value = 15

# Row 18
$ws.Range("D18").Value = (Get-Date -Year 2021 -Month 5 -Day 3 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M18").Value = 120

# Row 19
$ws.Range("D19").Value = (Get-Date -Year 2021 -Month 6 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("S19").Value = 1000

